$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 242, shifting rows 242:304 down to 243:305.
$ws.Rows.Item(242).Insert()

# Populate new row 242 with the new data entry.
$ws.Range("A242").Value = 4
$ws.Range("B242").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C242").Value = "Los Lagos"
$ws.Range("D242").Value = 44855
$ws.Range("E242").Value = 10
$ws.Range("F242").Value = 100112032
$ws.Range("G242").Value = "Zapallo italiano"
$ws.Range("H242").Value = "Sin especificar"
$ws.Range("I242").Value = "Primera"
$ws.Range("J242").Value = 200
$ws.Range("K242").Value = 18000
$ws.Range("L242").Value = 18000
$ws.Range("M242").Value = 18000
$ws.Range("N242").Value = "$/caja 50 unidades"
$ws.Range("O242").Value = "Región de Arica y Parinacota"
$ws.Range("P242").Value = 360
$ws.Range("Q242").Value = 50
$ws.Range("R242").Value = "Hortaliza"
